$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @(
    "2025-04-23 15:55:00:000",
    "2025-04-23 15:56:00:000",
    "2025-04-23 15:57:00:000",
    "2025-04-23 15:58:00:000",
    "2025-04-23 15:59:00:000",
    "2025-04-23 16:00:00:000"
)

$data = @(
    @(0.88949999999999996, 0.1951, -0.1986, 99, 0, 26.2, 174.36, 0.1265, 0.13450000000000001, 0.27639999999999998, 238),
    @(0.91159999999999997, 0.15479999999999999, 0.074800000000000005, 68, 0, 26, 204.3, 0.20039999999999999, 0.2482, 0.21829999999999999, 252),
    @(0.75329999999999997, 0.34410000000000002, -0.42980000000000002, 125, 0, 25.5, 206.34, 0.1651, 0.1026, 0.249, 223),
    @(0.75449999999999995, 0.38450000000000001, -0.44950000000000001, 143, 0, 25.5, 157.93, 0.1061, 0.13070000000000001, 0.10920000000000001, 282),
    @(0.57469999999999999, 0.069500000000000006, -0.61070000000000002, 103, 0, 25.5, 191.49, 0.35549999999999998, 0.1734, 0.3367, 223),
    @(0.52939999999999998, 0.1376, -0.3881, 175, 0, 25.2, 526.95000000000005, 0.45860000000000001, 0.52500000000000002, 0.33879999999999999, 371)
)

for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $timestamps[$i]

    $rowData = $data[$i]
    $ws.Cells.Item($row, 2).Value = $rowData[0]
    $ws.Cells.Item($row, 3).Value = $rowData[1]
    $ws.Cells.Item($row, 4).Value = $rowData[2]
    $ws.Cells.Item($row, 5).Value = $rowData[3]
    $ws.Cells.Item($row, 6).Value = $rowData[4]
    $ws.Cells.Item($row, 7).Value = $rowData[5]
    $ws.Cells.Item($row, 8).Value = $rowData[6]
    $ws.Cells.Item($row, 9).Value = $rowData[7]
    $ws.Cells.Item($row, 10).Value = $rowData[8]
    $ws.Cells.Item($row, 11).Value = $rowData[9]
    $ws.Cells.Item($row, 12).Value = $rowData[10]
}

$ws.Range("P14").Select()
